$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -0.1534545817258671
$ws.Range("D2").Value = 0.8794386829777929

$ws.Range("C3").Value = 1.858899280073114
$ws.Range("D3").Value = 0.07647067566880916

$ws.Range("C4").Value = 1.879444178685947
$ws.Range("D4").Value = 0.07349497010435901

$ws.Range("C5").Value = 3.886637124951239
$ws.Range("D5").Value = 0.0007948725780655419

$ws.Range("C6").Value = 2.013032339962128
$ws.Range("D6").Value = 0.05650739767483692
$ws.Range("G6").Value = "No"

$ws.Range("C7").Value = 1.983794753142436
$ws.Range("D7").Value = 0.05989484828259251

$ws.Range("C8").Value = 3.645023139327802
$ws.Range("D8").Value = 0.00142770096975231

$ws.Range("C9").Value = 0.1476646839790464
$ws.Range("D9").Value = 0.8839523907436035

$ws.Range("C10").Value = 2.742132000687696
$ws.Range("D10").Value = 0.01189581992882371

$ws.Range("C11").Value = 2.384827022983488
$ws.Range("D11").Value = 0.02613393036523814
